$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "2022-Q4" worksheet, right before the existing "2022-Q1"
#    sheet, and populate it with the new quarter's fund-holding detail.
# ---------------------------------------------------------------------------
$q1Sheet = $wb.Worksheets.Item(2)
$q4Sheet = $wb.Worksheets.Add($q1Sheet)
$q4Sheet.Name = "2022-Q4"

# NOTE: "$q1Sheet" now resolves to whatever sheet sits at position 2 (i.e.
# the brand-new "2022-Q4" sheet itself), since worksheet refs track tab
# position rather than the sheet they originally pointed to. Re-fetch the
# "2022-Q1" sheet by its new position (3) before using it as a copy source.
$q1Sheet = $wb.Worksheets.Item(3)

# Reuse the header row + row-2 layout/formatting from the "2022-Q1" sheet
# so the new sheet matches the existing look & feel (bold header, bordered
# first column, etc.)
$q1Sheet.Range("B1:H1").Copy($q4Sheet.Range("B1"))
$q1Sheet.Range("A2:H2").Copy($q4Sheet.Range("A2"))

# Fill in the actual 2022-Q4 holding for 迈得医疗 (688310)
$q4Sheet.Range("B2:G2").NumberFormat = "@"
$q4Sheet.Range("B2").Value = "007835"
$q4Sheet.Range("C2").Value = "国泰鑫睿混合"
$q4Sheet.Range("D2").Value = "7.62"
$q4Sheet.Range("E2").Value = "79.19"
$q4Sheet.Range("F2").Value = "2.96"
$q4Sheet.Range("G2").Value = "0.2256"
$q4Sheet.Range("H2").Value = 9

# ---------------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert a new row for 2022-Q4 above
#    the existing quarters, shifting 2022-Q1 / 2021-Q4 / 2021-Q3 down.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)
$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

# Copy row 3's A-column formatting down into the freshly inserted row 2 so
# the new row's index cell matches the existing style (bordered/centered).
$total.Range("A3").Copy($total.Range("A2"))

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.23
